# Auto-generated edit script
# Updates market-price-derived columns (H-N) across the Kujata_Profits sheets
# to reflect refreshed Universalis price data, per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = @(
    @{ Addr = "H17"; Value = 2680.6365 }
    @{ Addr = "J17"; Value = 2680.6365 }
    @{ Addr = "L17"; Value = 8041.9095 }
    @{ Addr = "N17"; Value = -8377.9095 }
    @{ Addr = "H69"; Value = 3366.875 }
    @{ Addr = "J69"; Value = 3419.2856 }
    @{ Addr = "L69"; Value = 10257.8568 }
    @{ Addr = "N69"; Value = -12005.8568 }
    @{ Addr = "H70"; Value = 1649.4445 }
    @{ Addr = "J70"; Value = 1600.75 }
    @{ Addr = "L70"; Value = 4802.25 }
    @{ Addr = "N70"; Value = -5342.25 }
    @{ Addr = "H72"; Value = 3366.875 }
    @{ Addr = "J72"; Value = 3419.2856 }
    @{ Addr = "L72"; Value = 30773.5704 }
    @{ Addr = "N72"; Value = -39509.5704 }
    @{ Addr = "H73"; Value = 1649.4445 }
    @{ Addr = "J73"; Value = 1600.75 }
    @{ Addr = "L73"; Value = 4802.25 }
    @{ Addr = "N73"; Value = -6674.25 }
    @{ Addr = "H86"; Value = 5140 }
    @{ Addr = "I86"; Value = 7233.3335 }
    @{ Addr = "J86"; Value = 2000 }
    @{ Addr = "K86"; Value = 7233.3335 }
    @{ Addr = "L86"; Value = 2000 }
    @{ Addr = "M86"; Value = -6110.3335 }
    @{ Addr = "N86"; Value = -4246 }
    @{ Addr = "H89"; Value = 5140 }
    @{ Addr = "I89"; Value = 7233.3335 }
    @{ Addr = "J89"; Value = 2000 }
    @{ Addr = "K89"; Value = 36166.6675 }
    @{ Addr = "L89"; Value = 10000 }
    @{ Addr = "M89"; Value = -30550.6675 }
    @{ Addr = "N89"; Value = -21232 }
    @{ Addr = "H137"; Value = 1777.0625 }
    @{ Addr = "I137"; Value = 1020.8 }
    @{ Addr = "J137"; Value = 2120.818 }
    @{ Addr = "K137"; Value = 3062.4 }
    @{ Addr = "L137"; Value = 6362.454000000001 }
    @{ Addr = "M137"; Value = -512.3999999999996 }
    @{ Addr = "N137"; Value = -11462.454 }
    @{ Addr = "H138"; Value = 545240.75 }
    @{ Addr = "I138"; Value = 1673.0454 }
    @{ Addr = "J138"; Value = 751421.6 }
    @{ Addr = "K138"; Value = 5019.1362 }
    @{ Addr = "L138"; Value = 2254264.8 }
    @{ Addr = "M138"; Value = 120.8638000000001 }
    @{ Addr = "N138"; Value = -2264544.8 }
)
foreach ($u in $ALC_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = @(
    @{ Addr = "H61"; Value = 1623.8572 }
    @{ Addr = "I61"; Value = 1652.5714 }
    @{ Addr = "J61"; Value = 1537.7142 }
    @{ Addr = "K61"; Value = 1652.5714 }
    @{ Addr = "L61"; Value = 1537.7142 }
    @{ Addr = "M61"; Value = -1440.5714 }
    @{ Addr = "N61"; Value = -1961.7142 }
    @{ Addr = "H97"; Value = 658.0625 }
    @{ Addr = "I97"; Value = 502.41666 }
    @{ Addr = "J97"; Value = 1125 }
    @{ Addr = "K97"; Value = 502.41666 }
    @{ Addr = "L97"; Value = 1125 }
    @{ Addr = "M97"; Value = -6.416659999999979 }
    @{ Addr = "N97"; Value = -2117 }
    @{ Addr = "H136"; Value = 1623.8572 }
    @{ Addr = "I136"; Value = 1652.5714 }
    @{ Addr = "J136"; Value = 1537.7142 }
    @{ Addr = "K136"; Value = 4957.7142 }
    @{ Addr = "L136"; Value = 4613.142599999999 }
    @{ Addr = "M136"; Value = -2407.7142 }
    @{ Addr = "N136"; Value = -9713.142599999999 }
)
foreach ($u in $ARM_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = @(
    @{ Addr = "H86"; Value = 5820.8 }
    @{ Addr = "I86"; Value = 5820.8 }
    @{ Addr = "K86"; Value = 5820.8 }
    @{ Addr = "M86"; Value = -4697.8 }
    @{ Addr = "H89"; Value = 5820.8 }
    @{ Addr = "I89"; Value = 5820.8 }
    @{ Addr = "K89"; Value = 29104 }
    @{ Addr = "M89"; Value = -23488 }
    @{ Addr = "H94"; Value = 13158683 }
    @{ Addr = "I94"; Value = 13889597 }
    @{ Addr = "K94"; Value = 13889597 }
    @{ Addr = "M94"; Value = -13889146 }
    @{ Addr = "H134"; Value = 5810.6787 }
    @{ Addr = "I134"; Value = 1329.25 }
    @{ Addr = "K134"; Value = 3987.75 }
    @{ Addr = "M134"; Value = -1452.75 }
)
foreach ($u in $BSM_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = @(
    @{ Addr = "H31"; Value = 1579.7778 }
    @{ Addr = "I31"; Value = 1572 }
    @{ Addr = "J31"; Value = 1607 }
    @{ Addr = "K31"; Value = 1572 }
    @{ Addr = "L31"; Value = 1607 }
    @{ Addr = "M31"; Value = -1277 }
    @{ Addr = "N31"; Value = -2197 }
    @{ Addr = "H34"; Value = 1579.7778 }
    @{ Addr = "I34"; Value = 1572 }
    @{ Addr = "J34"; Value = 1607 }
    @{ Addr = "K34"; Value = 1572 }
    @{ Addr = "L34"; Value = 1607 }
    @{ Addr = "M34"; Value = -1370 }
    @{ Addr = "N34"; Value = -2011 }
    @{ Addr = "H58"; Value = 549.2941 }
    @{ Addr = "I58"; Value = 548.5 }
    @{ Addr = "J58"; Value = 550.4286 }
    @{ Addr = "K58"; Value = 548.5 }
    @{ Addr = "L58"; Value = 550.4286 }
    @{ Addr = "M58"; Value = -345.5 }
    @{ Addr = "N58"; Value = -956.4286 }
    @{ Addr = "H132"; Value = 1673.8292 }
    @{ Addr = "I132"; Value = 1296.6364 }
    @{ Addr = "K132"; Value = 3889.9092 }
    @{ Addr = "M132"; Value = -1359.9092 }
    @{ Addr = "H134"; Value = 649.6842 }
    @{ Addr = "I134"; Value = 582.5143 }
    @{ Addr = "K134"; Value = 1747.5429 }
    @{ Addr = "M134"; Value = 787.4570999999999 }
    @{ Addr = "H136"; Value = 549.2941 }
    @{ Addr = "I136"; Value = 548.5 }
    @{ Addr = "J136"; Value = 550.4286 }
    @{ Addr = "K136"; Value = 1645.5 }
    @{ Addr = "L136"; Value = 1651.2858 }
    @{ Addr = "M136"; Value = 904.5 }
    @{ Addr = "N136"; Value = -6751.2858 }
    @{ Addr = "H141"; Value = 27560 }
    @{ Addr = "J141"; Value = 27560 }
    @{ Addr = "L141"; Value = 27560 }
    @{ Addr = "N141"; Value = -37920 }
)
foreach ($u in $CRP_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = @(
    @{ Addr = "H4"; Value = 683475.6 }
    @{ Addr = "J4"; Value = 1366801.2 }
    @{ Addr = "L4"; Value = 4100403.6 }
    @{ Addr = "N4"; Value = -4100627.6 }
    @{ Addr = "H107"; Value = 8149.923 }
    @{ Addr = "I107"; Value = 441 }
    @{ Addr = "K107"; Value = 1323 }
    @{ Addr = "M107"; Value = 597 }
    @{ Addr = "H119"; Value = 7158.923 }
    @{ Addr = "I119"; Value = 3016.5 }
    @{ Addr = "J119"; Value = 9000 }
    @{ Addr = "K119"; Value = 9049.5 }
    @{ Addr = "L119"; Value = 27000 }
    @{ Addr = "M119"; Value = -4211.5 }
    @{ Addr = "N119"; Value = -36676 }
    @{ Addr = "H120"; Value = 9574.5 }
    @{ Addr = "J120"; Value = 11999.333 }
    @{ Addr = "L120"; Value = 35997.999 }
    @{ Addr = "N120"; Value = -45673.999 }
    @{ Addr = "H122"; Value = 835.2941 }
    @{ Addr = "I122"; Value = 510 }
    @{ Addr = "J122"; Value = 1124.4445 }
    @{ Addr = "K122"; Value = 4590 }
    @{ Addr = "L122"; Value = 10120.0005 }
    @{ Addr = "M122"; Value = -2140 }
    @{ Addr = "N122"; Value = -15020.0005 }
)
foreach ($u in $CUL_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = @(
    @{ Addr = "H126"; Value = 2600 }
    @{ Addr = "H132"; Value = 2352.9565 }
    @{ Addr = "I132"; Value = 1487.8667 }
    @{ Addr = "K132"; Value = 4463.6001 }
    @{ Addr = "M132"; Value = -1933.6001 }
)
foreach ($u in $GSM_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = @(
    @{ Addr = "H16"; Value = 943.6923 }
    @{ Addr = "I16"; Value = 832.1 }
    @{ Addr = "J16"; Value = 1315.6666 }
    @{ Addr = "K16"; Value = 832.1 }
    @{ Addr = "L16"; Value = 1315.6666 }
    @{ Addr = "M16"; Value = -662.1 }
    @{ Addr = "N16"; Value = -1655.6666 }
    @{ Addr = "H55"; Value = 370.55 }
    @{ Addr = "I55"; Value = 278.15384 }
    @{ Addr = "K55"; Value = 278.15384 }
    @{ Addr = "M55"; Value = -105.15384 }
    @{ Addr = "H61"; Value = 2251.4167 }
    @{ Addr = "I61"; Value = 1847.5 }
    @{ Addr = "J61"; Value = 2453.375 }
    @{ Addr = "K61"; Value = 1847.5 }
    @{ Addr = "L61"; Value = 2453.375 }
    @{ Addr = "M61"; Value = -1645.5 }
    @{ Addr = "N61"; Value = -2857.375 }
    @{ Addr = "H68"; Value = 2067.2222 }
    @{ Addr = "I68"; Value = 1650.8334 }
    @{ Addr = "J68"; Value = 2900 }
    @{ Addr = "K68"; Value = 1650.8334 }
    @{ Addr = "L68"; Value = 2900 }
    @{ Addr = "M68"; Value = -901.8334 }
    @{ Addr = "N68"; Value = -4398 }
    @{ Addr = "H71"; Value = 2067.2222 }
    @{ Addr = "I71"; Value = 1650.8334 }
    @{ Addr = "J71"; Value = 2900 }
    @{ Addr = "K71"; Value = 8254.166999999999 }
    @{ Addr = "L71"; Value = 14500 }
    @{ Addr = "M71"; Value = -4510.166999999999 }
    @{ Addr = "N71"; Value = -21988 }
    @{ Addr = "H113"; Value = 2251.4167 }
    @{ Addr = "I113"; Value = 1847.5 }
    @{ Addr = "J113"; Value = 2453.375 }
    @{ Addr = "K113"; Value = 1847.5 }
    @{ Addr = "L113"; Value = 2453.375 }
    @{ Addr = "M113"; Value = 322.5 }
    @{ Addr = "N113"; Value = -6793.375 }
)
foreach ($u in $LTW_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = @(
    @{ Addr = "H132"; Value = 2191.7046 }
    @{ Addr = "I132"; Value = 2058.024 }
    @{ Addr = "K132"; Value = 6174.072 }
    @{ Addr = "M132"; Value = -3644.072 }
)
foreach ($u in $WVR_updates) {
    $ws.Range($u.Addr).Value = $u.Value
}
